# Re-sort the account-statement table: group rows by worker (grouped
# together instead of interleaved), each worker's periods in descending
# order. The underlying set of records is unchanged -- only the row
# order (and therefore which shared-string/value ends up in which row)
# changes.
#
# Commit message: "Elimna EC anteriores y se agregan nuevos, se modifica
# base de datos" -- previous EC rows removed, new ones added, data
# re-sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 16; Doc = "CC"; Id = "1045731943"; Name = "DALGIS MARIA LOBO LARA"; Period = "2112"; Mora = 18726;  Salario = 877803  },
    @{ Row = 17; Doc = "CC"; Id = "1045731943"; Name = "DALGIS MARIA LOBO LARA"; Period = "2111"; Mora = 35112;  Salario = 877803  },
    @{ Row = 18; Doc = "CC"; Id = "1045731943"; Name = "DALGIS MARIA LOBO LARA"; Period = "2110"; Mora = 35112;  Salario = 877803  },
    @{ Row = 19; Doc = "CC"; Id = "1045731943"; Name = "DALGIS MARIA LOBO LARA"; Period = "2109"; Mora = 35112;  Salario = 877803  },
    @{ Row = 20; Doc = "CC"; Id = "1045731943"; Name = "DALGIS MARIA LOBO LARA"; Period = "2108"; Mora = 35112;  Salario = 877803  },
    @{ Row = 21; Doc = "CC"; Id = "1045731943"; Name = "DALGIS MARIA LOBO LARA"; Period = "2105"; Mora = 35112;  Salario = 877803  },
    @{ Row = 22; Doc = "CC"; Id = "1045731943"; Name = "DALGIS MARIA LOBO LARA"; Period = "2104"; Mora = 35112;  Salario = 877803  },
    @{ Row = 23; Doc = "CC"; Id = "8373933";    Name = "ENDER MARTINEZ DIAZ";    Period = "2112"; Mora = 32000;  Salario = 1500000 },
    @{ Row = 24; Doc = "CC"; Id = "8373933";    Name = "ENDER MARTINEZ DIAZ";    Period = "2111"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 25; Doc = "CC"; Id = "8373933";    Name = "ENDER MARTINEZ DIAZ";    Period = "2110"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 26; Doc = "CC"; Id = "8373933";    Name = "ENDER MARTINEZ DIAZ";    Period = "2109"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 27; Doc = "CC"; Id = "8373933";    Name = "ENDER MARTINEZ DIAZ";    Period = "2108"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 28; Doc = "CC"; Id = "8373933";    Name = "ENDER MARTINEZ DIAZ";    Period = "2105"; Mora = 60000;  Salario = 1500000 },
    @{ Row = 29; Doc = "CC"; Id = "8373933";    Name = "ENDER MARTINEZ DIAZ";    Period = "2104"; Mora = 60000;  Salario = 1500000 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Range("B$r").Value = $rec.Doc
    $ws.Range("C$r").Value = $rec.Id
    $ws.Range("D$r").Value = $rec.Name
    $ws.Range("E$r").Value = $rec.Period
    $ws.Range("F$r").Value = $rec.Mora
    $ws.Range("G$r").Value = $rec.Salario
}

# Columns B/C/E/F/G are bestFit -- let Excel recompute the displayed
# column widths now that row 16's (and others') contents changed.
$ws.Range("B16:J29").Columns.AutoFit() | Out-Null
